$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Rename the existing sheet "Sheet1" -> "1.4C"
# ------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "1.4C"

# ------------------------------------------------------------------
# 2. Add a new sheet "1.4D" right after "1.4C" and populate it with
#    the same kind of car-manufacturing table, shifted one column to
#    the left (A..I instead of B..J) compared to "1.4C".
# ------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "1.4D"

# Headers (row 1) - bold, matches the style used on "1.4C"'s header row
$ws2.Range("A1").Value2 = "Car"
$ws2.Range("B1").Value2 = "time to manufacture"
$ws2.Range("C1").Value2 = "profit"
$ws2.Range("D1").Value2 = "orders"
$ws2.Range("E1").Value2 = "hours available"
$ws2.Range("G1").Value2 = "amount to manufacture"
$ws2.Range("H1").Value2 = "Total_profit"
$ws2.Range("I1").Value2 = "Hours to manufacture"
$ws2.Range("A1:I1").Font.Bold = $true

# Row 2 - Tesla
$ws2.Range("A2").Value2 = "Tesla"
$ws2.Range("B2").Value2 = 1
$ws2.Range("C2").Value2 = 200
$ws2.Range("D2").Value2 = 10
$ws2.Range("G2").Value2 = 34.999999442509001
$ws2.Range("H2").Formula = "=G2*C2"
$ws2.Range("I2").Formula = "=G2*B2"

# Row 3 - Chevy
$ws2.Range("A3").Value2 = "Chevy"
$ws2.Range("B3").Value2 = 2
$ws2.Range("C3").Value2 = 500
$ws2.Range("D3").Value2 = 20
$ws2.Range("G3").Value2 = 20
$ws2.Range("H3").Formula = "=G3*C3"
$ws2.Range("I3").Formula = "=G3*B3"

# Row 4 - Lexus
$ws2.Range("A4").Value2 = "Lexus"
$ws2.Range("B4").Value2 = 3
$ws2.Range("C4").Value2 = 700
$ws2.Range("D4").Value2 = 15
$ws2.Range("G4").Value2 = 15
$ws2.Range("H4").Formula = "=G4*C4"
$ws2.Range("I4").Formula = "=G4*B4"

# Row 5 - totals
$ws2.Range("E5").Value2 = 120
$ws2.Range("G5").Formula = "=SUM(G2:G4)"
$ws2.Range("H5").Formula = "=SUM(H2:H4)"
$ws2.Range("I5").Formula = "=SUM(I2:I4)"

# Column widths, matching "1.4C"'s corresponding (shifted) columns
$ws2.Columns.Item(2).ColumnWidth = $ws1.Columns.Item(3).ColumnWidth
$ws2.Columns.Item(7).ColumnWidth = $ws1.Columns.Item(8).ColumnWidth
$ws2.Columns.Item(9).ColumnWidth = $ws1.Columns.Item(10).ColumnWidth

# ------------------------------------------------------------------
# 3. Solver parameters (hidden defined names) for each sheet
# ------------------------------------------------------------------
$adj1 = $ws1.Names.Add("solver_adj", "='1.4C'!`$H`$2:`$H`$4")
$adj1.Visible = $false
$cvg1 = $ws1.Names.Add("solver_cvg", "=0.0001")
$cvg1.Visible = $false
$drv1 = $ws1.Names.Add("solver_drv", "=1")
$drv1.Visible = $false
$eng1 = $ws1.Names.Add("solver_eng", "=1")
$eng1.Visible = $false
$itr1 = $ws1.Names.Add("solver_itr", "=2147483647")
$itr1.Visible = $false
$lhs11 = $ws1.Names.Add("solver_lhs1", "='1.4C'!`$H`$2:`$H`$4")
$lhs11.Visible = $false
$lhs21 = $ws1.Names.Add("solver_lhs2", "='1.4C'!`$J`$5")
$lhs21.Visible = $false
$lin1 = $ws1.Names.Add("solver_lin", "=2")
$lin1.Visible = $false
$mip1 = $ws1.Names.Add("solver_mip", "=2147483647")
$mip1.Visible = $false
$mni1 = $ws1.Names.Add("solver_mni", "=30")
$mni1.Visible = $false
$mrt1 = $ws1.Names.Add("solver_mrt", "=0.075")
$mrt1.Visible = $false
$msl1 = $ws1.Names.Add("solver_msl", "=2")
$msl1.Visible = $false
$neg1 = $ws1.Names.Add("solver_neg", "=1")
$neg1.Visible = $false
$nod1 = $ws1.Names.Add("solver_nod", "=2147483647")
$nod1.Visible = $false
$num1 = $ws1.Names.Add("solver_num", "=2")
$num1.Visible = $false
$opt1 = $ws1.Names.Add("solver_opt", "='1.4C'!`$I`$5")
$opt1.Visible = $false
$pre1 = $ws1.Names.Add("solver_pre", "=0.000001")
$pre1.Visible = $false
$rbv1 = $ws1.Names.Add("solver_rbv", "=1")
$rbv1.Visible = $false
$rel11 = $ws1.Names.Add("solver_rel1", "=3")
$rel11.Visible = $false
$rel21 = $ws1.Names.Add("solver_rel2", "=1")
$rel21.Visible = $false
$rhs11 = $ws1.Names.Add("solver_rhs1", "='1.4C'!`$E`$2:`$E`$4")
$rhs11.Visible = $false
$rhs21 = $ws1.Names.Add("solver_rhs2", "='1.4C'!`$F`$5")
$rhs21.Visible = $false
$rlx1 = $ws1.Names.Add("solver_rlx", "=2")
$rlx1.Visible = $false
$rsd1 = $ws1.Names.Add("solver_rsd", "=0")
$rsd1.Visible = $false
$scl1 = $ws1.Names.Add("solver_scl", "=1")
$scl1.Visible = $false
$sho1 = $ws1.Names.Add("solver_sho", "=2")
$sho1.Visible = $false
$ssz1 = $ws1.Names.Add("solver_ssz", "=100")
$ssz1.Visible = $false
$tim1 = $ws1.Names.Add("solver_tim", "=2147483647")
$tim1.Visible = $false
$tol1 = $ws1.Names.Add("solver_tol", "=0.01")
$tol1.Visible = $false
$typ1 = $ws1.Names.Add("solver_typ", "=1")
$typ1.Visible = $false
$val1 = $ws1.Names.Add("solver_val", "=0")
$val1.Visible = $false
$ver1 = $ws1.Names.Add("solver_ver", "=2")
$ver1.Visible = $false

$adj2 = $ws2.Names.Add("solver_adj", "='1.4D'!`$G`$2:`$G`$4")
$adj2.Visible = $false
$cvg2 = $ws2.Names.Add("solver_cvg", "=0.0001")
$cvg2.Visible = $false
$drv2 = $ws2.Names.Add("solver_drv", "=1")
$drv2.Visible = $false
$eng2 = $ws2.Names.Add("solver_eng", "=1")
$eng2.Visible = $false
$itr2 = $ws2.Names.Add("solver_itr", "=2147483647")
$itr2.Visible = $false
$lhs12 = $ws2.Names.Add("solver_lhs1", "='1.4D'!`$G`$2:`$G`$4")
$lhs12.Visible = $false
$lhs22 = $ws2.Names.Add("solver_lhs2", "='1.4D'!`$I`$5")
$lhs22.Visible = $false
$lin2 = $ws2.Names.Add("solver_lin", "=2")
$lin2.Visible = $false
$mip2 = $ws2.Names.Add("solver_mip", "=2147483647")
$mip2.Visible = $false
$mni2 = $ws2.Names.Add("solver_mni", "=30")
$mni2.Visible = $false
$mrt2 = $ws2.Names.Add("solver_mrt", "=0.075")
$mrt2.Visible = $false
$msl2 = $ws2.Names.Add("solver_msl", "=2")
$msl2.Visible = $false
$neg2 = $ws2.Names.Add("solver_neg", "=1")
$neg2.Visible = $false
$nod2 = $ws2.Names.Add("solver_nod", "=2147483647")
$nod2.Visible = $false
$num2 = $ws2.Names.Add("solver_num", "=2")
$num2.Visible = $false
$opt2 = $ws2.Names.Add("solver_opt", "='1.4D'!`$G`$5")
$opt2.Visible = $false
$pre2 = $ws2.Names.Add("solver_pre", "=0.000001")
$pre2.Visible = $false
$rbv2 = $ws2.Names.Add("solver_rbv", "=1")
$rbv2.Visible = $false
$rel12 = $ws2.Names.Add("solver_rel1", "=3")
$rel12.Visible = $false
$rel22 = $ws2.Names.Add("solver_rel2", "=1")
$rel22.Visible = $false
$rhs12 = $ws2.Names.Add("solver_rhs1", "='1.4D'!`$D`$2:`$D`$4")
$rhs12.Visible = $false
$rhs22 = $ws2.Names.Add("solver_rhs2", "='1.4D'!`$E`$5")
$rhs22.Visible = $false
$rlx2 = $ws2.Names.Add("solver_rlx", "=2")
$rlx2.Visible = $false
$rsd2 = $ws2.Names.Add("solver_rsd", "=0")
$rsd2.Visible = $false
$scl2 = $ws2.Names.Add("solver_scl", "=1")
$scl2.Visible = $false
$sho2 = $ws2.Names.Add("solver_sho", "=2")
$sho2.Visible = $false
$ssz2 = $ws2.Names.Add("solver_ssz", "=100")
$ssz2.Visible = $false
$tim2 = $ws2.Names.Add("solver_tim", "=2147483647")
$tim2.Visible = $false
$tol2 = $ws2.Names.Add("solver_tol", "=0.01")
$tol2.Visible = $false
$typ2 = $ws2.Names.Add("solver_typ", "=1")
$typ2.Visible = $false
$val2 = $ws2.Names.Add("solver_val", "=0")
$val2.Visible = $false
$ver2 = $ws2.Names.Add("solver_ver", "=2")
$ver2.Visible = $false

# ------------------------------------------------------------------
# 4. View state: selections + zoom + active sheet
# ------------------------------------------------------------------
$ws2.Range("D6").Select()
$excel.ActiveWindow.Zoom = 161

$ws1.Activate()
$ws1.Range("E15").Select()

$ws2.Activate()
